# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45206 to 45208, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 103
$col = 3  # Column C = "Förändrad"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $col).Value = 45208
}
